$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "30.319.87"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  +2.12%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.089.11"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  -0.27%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.002"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  -0.77%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "343.19"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  -0.55%  "
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  -0.68%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5237"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  +1.78%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.4421"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  +0.71%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "54.57"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  +3.86%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.09328"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  +0.72%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.169"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  -0.16%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "24.86"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  +0.07%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "8.602"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  +3.88%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.897"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  +2.21%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "2.088.44"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  -0.41%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "101.26"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  +1.77%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.00001160"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  +0.67%  "
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  -0.70%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "21.13"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  +1.24%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.06664"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  +0.12%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.337"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  +2.28%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "1.002"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  -0.57%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "30.309.02"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  +1.92%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "12.54"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  -0.34%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.303"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  -0.67%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "21.80"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  -0.50%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "162.69"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  +0.48%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.517"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  -0.38%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "133.13"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  +0.06%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.138"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  +0.59%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.668"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  +1.08%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.1045"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  -0.50%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "6.240"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  +1.17%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "6.708"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  +8.72%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "3.856"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  -2.04%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "10.17"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  -1.36%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.02628"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  +2.22%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.06846"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  +2.11%  "
$ws.Range("B39").Value = "TheSandbox"
$ws.Range("C39").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.6986"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  +1.92%  "
$ws.Range("B40").Value = "TrustWalletToken"
$ws.Range("C40").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.347"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  +3.94%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "12.52"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  +0.61%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.2214"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  -0.39%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.6825"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  +2.92%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "14.36"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  +0.00%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.339"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  +1.01%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.001"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  -0.62%  "
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  +18.40%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "3.634"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  +0.19%  "
$ws.Range("B49").Value = "ThetaToken"
$ws.Range("C49").Value = "https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.210"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  +8.07%  "
$ws.Range("B50").Value = "BabyDogeCoin"
$ws.Range("C50").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.00000000342"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  -1.48%  "
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  -0.37%  "
